$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "TERRACED"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "TERRACED"

$ws.Range("D9").Select()
